# Update the cryptos worksheet with freshly scraped price/volume figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 and 17 swap coin identities (ShibaInu <-> WrappedEther) in addition
# to getting new price/volume figures, so set B/C/D/E fully for those rows.
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.00001026"
$ws.Range("E16").Value = "  -1.42%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "1.422.22"
$ws.Range("E17").Value = "  -4.50%  "

# All other rows only change the Price (D) and Volume(1h) (E) columns.
$ws.Range("D2").Value = "20.011.40"
$ws.Range("E2").Value = "  -4.05%  "

$ws.Range("D3").Value = "1.421.39"
$ws.Range("E3").Value = "  -4.38%  "

$ws.Range("D4").Value = "'0.9995"
$ws.Range("E4").Value = "  -0.66%  "

$ws.Range("D5").Value = "'0.9996"
$ws.Range("E5").Value = "  -0.79%  "

$ws.Range("D6").Value = "'276.98"
$ws.Range("E6").Value = "  -1.14%  "

$ws.Range("D7").Value = "'0.3694"
$ws.Range("E7").Value = "  -2.44%  "

$ws.Range("D8").Value = "'0.3108"
$ws.Range("E8").Value = "  +0.61%  "

$ws.Range("D9").Value = "'39.85"

$ws.Range("D10").Value = "'1.050"
$ws.Range("E10").Value = "  +2.62%  "

$ws.Range("D11").Value = "'0.06557"
$ws.Range("E11").Value = "  -2.95%  "

$ws.Range("D12").Value = "'0.9993"
$ws.Range("E12").Value = "  -0.73%  "

$ws.Range("D13").Value = "'5.524"
$ws.Range("E13").Value = "  +0.78%  "

$ws.Range("D14").Value = "'17.78"
$ws.Range("E14").Value = "  +0.40%  "

$ws.Range("D15").Value = "'6.234"
$ws.Range("E15").Value = "  -0.59%  "

$ws.Range("D18").Value = "'0.05698"
$ws.Range("E18").Value = "  -12.46%  "

$ws.Range("D19").Value = "'0.9993"
$ws.Range("E19").Value = "  -0.87%  "

$ws.Range("D20").Value = "'71.57"
$ws.Range("E20").Value = "  -10.83%  "

$ws.Range("D21").Value = "'5.632"
$ws.Range("E21").Value = "  -4.31%  "

$ws.Range("D22").Value = "'14.82"
$ws.Range("E22").Value = "  -0.03%  "

$ws.Range("D23").Value = "'11.03"
$ws.Range("E23").Value = "  +2.19%  "

$ws.Range("D24").Value = "'2.241"
$ws.Range("E24").Value = "  -3.75%  "

$ws.Range("D25").Value = "20.029.25"
$ws.Range("E25").Value = "  -3.90%  "

$ws.Range("D26").Value = "'2.296"
$ws.Range("E26").Value = "  +0.19%  "

$ws.Range("D27").Value = "'133.34"
$ws.Range("E27").Value = "  -8.34%  "

$ws.Range("D28").Value = "'17.40"
$ws.Range("E28").Value = "  -1.88%  "

$ws.Range("D29").Value = "1.581.22"
$ws.Range("E29").Value = "  -4.63%  "

$ws.Range("D30").Value = "'110.41"
$ws.Range("E30").Value = "  -1.82%  "

$ws.Range("D31").Value = "'3.913"
$ws.Range("E31").Value = "  -18.36%  "

$ws.Range("D32").Value = "'5.273"
$ws.Range("E32").Value = "  -8.20%  "

$ws.Range("D33").Value = "'0.8254"
$ws.Range("E33").Value = "  -10.37%  "

$ws.Range("D34").Value = "'0.07779"
$ws.Range("E34").Value = "  -1.21%  "

$ws.Range("D35").Value = "'1.484"
$ws.Range("E35").Value = "  +2.09%  "

$ws.Range("D36").Value = "'8.248"
$ws.Range("E36").Value = "  -1.76%  "

$ws.Range("E37").Value = "  +1.02%  "

$ws.Range("D38").Value = "'0.05878"
$ws.Range("E38").Value = "  +2.36%  "

$ws.Range("D39").Value = "'0.9992"
$ws.Range("E39").Value = "  -0.74%  "

$ws.Range("D40").Value = "'0.02065"
$ws.Range("E40").Value = "  -0.25%  "

$ws.Range("D41").Value = "'10.55"
$ws.Range("E41").Value = "  -5.57%  "

$ws.Range("D42").Value = "'0.1891"
$ws.Range("E42").Value = "  -2.82%  "

$ws.Range("D43").Value = "'1.105"
$ws.Range("E43").Value = "  -2.72%  "

$ws.Range("D44").Value = "'0.5342"
$ws.Range("E44").Value = "  -2.74%  "

$ws.Range("D45").Value = "'12.43"
$ws.Range("E45").Value = "  -0.88%  "

$ws.Range("D46").Value = "'3.548"
$ws.Range("E46").Value = "  -2.27%  "

$ws.Range("D47").Value = "'117.27"
$ws.Range("E47").Value = "  +5.15%  "

$ws.Range("D48").Value = "'0.5225"
$ws.Range("E48").Value = "  -1.29%  "

$ws.Range("D49").Value = "'1.781"
$ws.Range("E49").Value = "  -2.08%  "

$ws.Range("E50").Value = "  -6.21%  "

$ws.Range("D51").Value = "'0.9991"
$ws.Range("E51").Value = "  -0.67%  "
